$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = "Normal"
}

Set-TextValue 'D2' '23.733.00'
Set-TextValue 'E2' '  +1.31%  '
Set-TextValue 'D3' '1.654.20'
Set-TextValue 'E3' '  +1.13%  '
Set-TextValue 'E4' '  +0.05%  '
Set-TextValue 'E5' '  +0.06%  '
Set-TextValue 'D6' '303.31'
Set-TextValue 'E6' '  -0.09%  '
Set-TextValue 'D7' '0.3805'
Set-TextValue 'E7' '  +0.58%  '
Set-TextValue 'D8' '0.3625'
Set-TextValue 'E8' '  -0.07%  '
Set-TextValue 'D9' '51.09'
Set-TextValue 'E9' '  -1.47%  '
Set-TextValue 'D10' '1.250'
Set-TextValue 'E10' '  +1.75%  '
Set-TextValue 'D11' '0.08216'
Set-TextValue 'E11' '  +0.45%  '
Set-TextValue 'D12' '1.001'
Set-TextValue 'E12' '  +0.05%  '
Set-TextValue 'D13' '22.70'
Set-TextValue 'E13' '  +1.00%  '
Set-TextValue 'D14' '6.518'
Set-TextValue 'E14' '  +0.74%  '
Set-TextValue 'D15' '7.441'
Set-TextValue 'E15' '  +0.65%  '
Set-TextValue 'D16' '0.00001235'
Set-TextValue 'E16' '  -0.37%  '
Set-TextValue 'D17' '1.653.59'
Set-TextValue 'E17' '  +1.51%  '
Set-TextValue 'D18' '97.40'
Set-TextValue 'E18' '  +2.54%  '
Set-TextValue 'D19' '0.07012'
Set-TextValue 'E19' '  +0.89%  '
Set-TextValue 'D20' '6.803'
Set-TextValue 'E20' '  +3.50%  '
Set-TextValue 'D21' '17.72'
Set-TextValue 'E21' '  +1.34%  '
Set-TextValue 'E22' '  +0.06%  '
Set-TextValue 'D23' '12.89'
Set-TextValue 'E23' '  +2.62%  '
Set-TextValue 'D24' '23.731.89'
Set-TextValue 'E24' '  +1.32%  '
Set-TextValue 'D25' '2.530'
Set-TextValue 'E25' '  +0.88%  '
Set-TextValue 'D26' '3.055'
Set-TextValue 'E26' '  +0.30%  '
Set-TextValue 'D27' '21.26'
Set-TextValue 'E27' '  +0.63%  '
Set-TextValue 'D28' '152.18'
Set-TextValue 'E28' '  +1.08%  '
Set-TextValue 'D29' '5.239'
Set-TextValue 'E29' '  -0.68%  '
Set-TextValue 'D30' '134.43'
Set-TextValue 'E30' '  +0.97%  '
Set-TextValue 'D31' '1.838.15'
Set-TextValue 'D32' '6.936'
Set-TextValue 'E32' '  +4.65%  '
Set-TextValue 'D33' '2.199'
Set-TextValue 'E33' '  +1.68%  '
Set-TextValue 'D34' '1.073'
Set-TextValue 'E34' '  +2.63%  '
Set-TextValue 'D35' '11.77'
Set-TextValue 'E35' '  +4.76%  '
Set-TextValue 'D36' '0.02812'
Set-TextValue 'E36' '  +2.14%  '
Set-TextValue 'D37' '0.2528'
Set-TextValue 'E37' '  +1.41%  '
Set-TextValue 'D38' '0.08815'
Set-TextValue 'E38' '  +0.42%  '
Set-TextValue 'E39' '  +1.34%  '
Set-TextValue 'D40' '0.07125'
Set-TextValue 'E40' '  +0.29%  '
Set-TextValue 'D41' '12.99'
Set-TextValue 'E41' '  +6.99%  '
Set-TextValue 'D42' '0.7039'
Set-TextValue 'E42' '  +0.52%  '
Set-TextValue 'D43' '1.339'
Set-TextValue 'E43' '  -0.27%  '
Set-TextValue 'D44' '16.04'
Set-TextValue 'E44' '  +1.61%  '
Set-TextValue 'D45' '0.6510'
Set-TextValue 'E45' '  +0.17%  '
Set-TextValue 'E46' '  +2.08%  '
Set-TextValue 'E47' '  +0.05%  '
Set-TextValue 'D48' '3.960'
Set-TextValue 'E48' '  -0.16%  '
Set-TextValue 'D49' '0.07960'
Set-TextValue 'E49' '  -0.20%  '
Set-TextValue 'D50' '128.48'
Set-TextValue 'E50' '  +1.27%  '
Set-TextValue 'D51' '1.192'
Set-TextValue 'E51' '  +0.27%  '
